$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed columns B (Fn1) and C (Itga2) and F=1 remain constant for all data rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "Fn1"
    $ws.Cells.Item($r, 3).Value = "Itga2"
    $ws.Cells.Item($r, 6).Value = 1
}

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 27.03890566666666
$ws.Cells.Item(2, 8).Value = 81.116717
$ws.Cells.Item(2, 9).Value = 0.07096188219033728
$ws.Cells.Item(2, 10).Value = 0.07096188219033729
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 1.888921
$ws.Cells.Item(2, 14).Value = 5.666763
$ws.Cells.Item(2, 15).Value = 0.4551532417350329
$ws.Cells.Item(2, 16).Value = 0.4551532417350328
$ws.Cells.Item(2, 17).Value = 51.07435673078566
$ws.Cells.Item(2, 18).Value = 459.6692105770709
$ws.Cells.Item(2, 19).Value = 0.03229853071855151
$ws.Cells.Item(2, 20).Value = 0.03229853071855151

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 27.03890566666666
$ws.Cells.Item(3, 8).Value = 81.116717
$ws.Cells.Item(3, 9).Value = 0.07096188219033728
$ws.Cells.Item(3, 10).Value = 0.07096188219033729
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.087098333333333
$ws.Cells.Item(3, 14).Value = 3.261295
$ws.Cells.Item(3, 15).Value = 0.2619465454094788
$ws.Cells.Item(3, 16).Value = 0.2619465454094787
$ws.Cells.Item(3, 17).Value = 29.39394928539055
$ws.Cells.Item(3, 18).Value = 264.545543568515
$ws.Cells.Item(3, 19).Value = 0.01858821989551327
$ws.Cells.Item(3, 20).Value = 0.01858821989551326

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 27.03890566666666
$ws.Cells.Item(4, 8).Value = 81.116717
$ws.Cells.Item(4, 9).Value = 0.07096188219033728
$ws.Cells.Item(4, 10).Value = 0.07096188219033729
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.174057666666667
$ws.Cells.Item(4, 14).Value = 3.522173
$ws.Cells.Item(4, 15).Value = 0.2829002128554884
$ws.Cells.Item(4, 16).Value = 0.2829002128554884
$ws.Cells.Item(4, 17).Value = 31.74523449622678
$ws.Cells.Item(4, 18).Value = 285.707110466041
$ws.Cells.Item(4, 19).Value = 0.0200751315762725
$ws.Cells.Item(4, 20).Value = 0.02007513157627251

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 345.566579
$ws.Cells.Item(5, 8).Value = 1036.699737
$ws.Cells.Item(5, 9).Value = 0.9069174311350353
$ws.Cells.Item(5, 10).Value = 0.9069174311350354
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 1.888921
$ws.Cells.Item(5, 14).Value = 5.666763
$ws.Cells.Item(5, 15).Value = 0.4551532417350329
$ws.Cells.Item(5, 16).Value = 0.4551532417350328
$ws.Cells.Item(5, 17).Value = 652.747967971259
$ws.Cells.Item(5, 18).Value = 5874.73171174133
$ws.Cells.Item(5, 19).Value = 0.4127864087671198
$ws.Cells.Item(5, 20).Value = 0.4127864087671198

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 345.566579
$ws.Cells.Item(6, 8).Value = 1036.699737
$ws.Cells.Item(6, 9).Value = 0.9069174311350353
$ws.Cells.Item(6, 10).Value = 0.9069174311350354
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.087098333333333
$ws.Cells.Item(6, 14).Value = 3.261295
$ws.Cells.Item(6, 15).Value = 0.2619465454094788
$ws.Cells.Item(6, 16).Value = 0.2619465454094787
$ws.Cells.Item(6, 17).Value = 375.6648520866016
$ws.Cells.Item(6, 18).Value = 3380.983668779415
$ws.Cells.Item(6, 19).Value = 0.2375638880574613
$ws.Cells.Item(6, 20).Value = 0.2375638880574613

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 345.566579
$ws.Cells.Item(7, 8).Value = 1036.699737
$ws.Cells.Item(7, 9).Value = 0.9069174311350353
$ws.Cells.Item(7, 10).Value = 0.9069174311350354
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.174057666666667
$ws.Cells.Item(7, 14).Value = 3.522173
$ws.Cells.Item(7, 15).Value = 0.2829002128554884
$ws.Cells.Item(7, 16).Value = 0.2829002128554884
$ws.Cells.Item(7, 17).Value = 405.7150914187223
$ws.Cells.Item(7, 18).Value = 3651.435822768501
$ws.Cells.Item(7, 19).Value = 0.2565671343104542
$ws.Cells.Item(7, 20).Value = 0.2565671343104542

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 8.428738666666668
$ws.Cells.Item(8, 8).Value = 25.286216
$ws.Cells.Item(8, 9).Value = 0.0221206866746274
$ws.Cells.Item(8, 10).Value = 0.02212068667462741
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 1.888921
$ws.Cells.Item(8, 14).Value = 5.666763
$ws.Cells.Item(8, 15).Value = 0.4551532417350329
$ws.Cells.Item(8, 16).Value = 0.4551532417350328
$ws.Cells.Item(8, 17).Value = 15.92122147097867
$ws.Cells.Item(8, 18).Value = 143.290993238808
$ws.Cells.Item(8, 19).Value = 0.01006830224936161
$ws.Cells.Item(8, 20).Value = 0.01006830224936161

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 8.428738666666668
$ws.Cells.Item(9, 8).Value = 25.286216
$ws.Cells.Item(9, 9).Value = 0.0221206866746274
$ws.Cells.Item(9, 10).Value = 0.02212068667462741
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.087098333333333
$ws.Cells.Item(9, 14).Value = 3.261295
$ws.Cells.Item(9, 15).Value = 0.2619465454094788
$ws.Cells.Item(9, 16).Value = 0.2619465454094787
$ws.Cells.Item(9, 17).Value = 9.162867756635556
$ws.Cells.Item(9, 18).Value = 82.46580980972
$ws.Cells.Item(9, 19).Value = 0.005794437456504138
$ws.Cells.Item(9, 20).Value = 0.005794437456504139

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 8.428738666666668
$ws.Cells.Item(10, 8).Value = 25.286216
$ws.Cells.Item(10, 9).Value = 0.0221206866746274
$ws.Cells.Item(10, 10).Value = 0.02212068667462741
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.174057666666667
$ws.Cells.Item(10, 14).Value = 3.522173
$ws.Cells.Item(10, 15).Value = 0.2829002128554884
$ws.Cells.Item(10, 16).Value = 0.2829002128554884
$ws.Cells.Item(10, 17).Value = 9.895825251929779
$ws.Cells.Item(10, 18).Value = 89.06242726736801
$ws.Cells.Item(10, 19).Value = 0.006257946968761657
$ws.Cells.Item(10, 20).Value = 0.006257946968761659
